$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price ticker refresh: Price (D) and Volume(1h) (E) columns.
# D-column values are plain text (e.g. "62.778.19", thousands-dot style)
# that Excel would otherwise auto-coerce to a number on assignment, so
# each is forced to Text format first and the cell style is reset to
# "Normal" afterwards so no stray style index is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.778.19'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.444.65'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.78%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.72'
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.444.12'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.41%  '
$ws.Range("E10").Value = '  +2.56%  '
$ws.Range("E11").Value = '  +2.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.28'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.62%  '
$ws.Range("E13").Value = '  +2.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.19%  '
$ws.Range("E15").Value = '  +5.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.889.56'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.84%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.635.84'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.440.59'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.37%  '
$ws.Range("E19").Value = '  -1.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '330.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.84%  '
$ws.Range("E22").Value = '  +1.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.05'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.46%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("E25").Value = '  +1.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '648.65'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +11.11%  '
$ws.Range("E27").Value = '  +17.85%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.55'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0989'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.565.63'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₆0492'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +72.72%  '
$ws.Range("E32").Value = '  +1.94%  '
$ws.Range("E33").Value = '  +6.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.87'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.76%  '
$ws.Range("E35").Value = '  +4.28%  '
$ws.Range("E36").Value = '  +1.23%  '
$ws.Range("E37").Value = '  +0.15%  '
$ws.Range("E38").Value = '  +3.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.51'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '153.38'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.16%  '
$ws.Range("E41").Value = '  +0.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '18.77'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.52%  '
$ws.Range("E43").Value = '  +8.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.76'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.50'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.99%  '
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("E47").Value = '  +27.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '145.16'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.50%  '
$ws.Range("E49").Value = '  +3.20%  '
$ws.Range("E50").Value = '  +5.56%  '
$ws.Range("E51").Value = '  +2.14%  '
